# Apply the "Shuffled Image data" update to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: new trial block header row (numeric values)
$ws.Range("A7").Value = 20230921
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = 2
$ws.Range("D7").Value = 3
$ws.Range("E7").Value = 4
$ws.Range("F7").Value = 5
$ws.Range("G7").Value = 6

# Row 8: grasp-frame lists for the 20230921 trial
$ws.Range("C8").Value = "4,5"
$ws.Range("D8").Value = "10,13,14"
$ws.Range("E8").Value = "14,18"
$ws.Range("F8").Value = "9,20"

# Row 9: new trial block header row (numeric values)
$ws.Range("A9").Value = 20230922
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = 2
$ws.Range("D9").Value = 3
$ws.Range("E9").Value = 4
$ws.Range("F9").Value = 5
$ws.Range("G9").Value = 6

# Row 10: grasp-frame lists / values for the 20230922 trial
$ws.Range("B10").Value = "14,18,19,20"
$ws.Range("C10").Value = "17,18,19,20"
$ws.Range("D10").Value = 3
$ws.Range("G10").Value = 13

# Row 11: new trial block header row (numeric values)
$ws.Range("A11").Value = 20231016
$ws.Range("B11").Value = 1
$ws.Range("D11").Value = 2
$ws.Range("F11").Value = 3

# Row 13: new trial block header row (numeric values)
$ws.Range("A13").Value = 20231101
$ws.Range("B13").Value = 1
$ws.Range("C13").Value = 3
$ws.Range("D13").Value = 4
$ws.Range("E13").Value = 6
$ws.Range("F13").Value = 2
$ws.Range("G13").Value = 5

# Row 14: grasp-frame lists for the 20231101 trial
# (entry order below intentionally matches the original authoring order so
# that new shared-string table entries are created in the same sequence)
$ws.Range("A14").Value = "1,3,1,2,3,2"
$ws.Range("B14").Value = "13,16,17,18,19"
$ws.Range("F14").Value = "8,9,18"
$ws.Range("C14").Value = "19,20"
$ws.Range("D14").Value = "9,18,20"
$ws.Range("G14").Value = "18,19"
$ws.Range("E14").Value = "5,16,17"

# Row 15: new trial block header row (numeric values)
$ws.Range("A15").Value = 20231103
$ws.Range("B15").Value = 2
$ws.Range("C15").Value = 9
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 8
$ws.Range("F15").Value = 5
$ws.Range("G15").Value = 10

# Row 16: grasp-frame lists for the 20231103 trial
$ws.Range("A16").Value = "2,1,3,2,1,3"
$ws.Range("B16").Value = "3,18"
$ws.Range("C16").Value = "12,13,14"
$ws.Range("E16").Value = "7,18"
$ws.Range("F16").Value = "9,20"
$ws.Range("G16").Value = 7

# Row 17: new trial block header row (numeric values)
$ws.Range("A17").Value = 20231201
$ws.Range("B17").Value = 3
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 4
$ws.Range("F17").Value = 2
$ws.Range("G17").Value = 5

# Row 18: grasp-frame lists for the 20231201 trial
$ws.Range("A18").Value = "2,3,1,2,3,1"
$ws.Range("B18").Value = "3,7,12,13,17,18,19"
$ws.Range("C18").Value = "2,4,5,8,9,10,15,16"
$ws.Range("D18").Value = "17,18,19,20"
$ws.Range("E18").Value = "13,19,20"
$ws.Range("F18").Value = "6,13,14,19"
$ws.Range("G18").Value = "7,13"

# Update the active selection / scroll position to match the author's final view
$ws.Range("G18").Select()
$excel.ActiveWindow.ScrollColumn = 4
$excel.ActiveWindow.ScrollRow = 1
